$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.121.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.39%  "

$ws.Range("D3").Value = "'1.909.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.47%  "

$ws.Range("D4").Value = "'0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'252.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "

$ws.Range("D6").Value = "'0.9987"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'0.5097"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.55%  "

$ws.Range("D8").Value = "'45.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.78%  "

$ws.Range("D9").Value = "'0.3031"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.74%  "

$ws.Range("D10").Value = "'0.06809"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.11%  "

$ws.Range("D11").Value = "'1.906.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.34%  "

$ws.Range("D12").Value = "'17.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.90%  "

$ws.Range("D13").Value = "'0.07322"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.38%  "

$ws.Range("D14").Value = "'0.6932"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.15%  "

$ws.Range("D15").Value = "'86.82"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'4.913"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.54%  "

$ws.Range("D17").Value = "'30.115.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.32%  "

$ws.Range("D18").Value = "'0.000008199"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.76%  "

$ws.Range("D19").Value = "'0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").Value = "'13.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.72%  "

$ws.Range("D21").Value = "'2.153.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.17%  "

$ws.Range("D22").Value = "'0.9979"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").Value = "'4.824"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.06%  "

$ws.Range("D24").Value = "'5.745"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.46%  "

$ws.Range("D25").Value = "'9.275"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.48%  "

$ws.Range("D26").Value = "'147.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.94%  "

$ws.Range("D27").Value = "'135.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.35%  "

$ws.Range("D28").Value = "'17.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.31%  "

$ws.Range("D29").Value = "'2.001"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.03%  "

$ws.Range("E30").Value = "  -0.69%  "

$ws.Range("D31").Value = "'4.270"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.08%  "

$ws.Range("D32").Value = "'0.08833"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.61%  "

$ws.Range("D33").Value = "'4.010"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.19%  "

$ws.Range("D34").Value = "'0.05054"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").Value = "'1.140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.99%  "

$ws.Range("D36").Value = "'0.7226"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.23%  "

$ws.Range("D37").Value = "'2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'2.821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.60%  "

$ws.Range("D39").Value = "'2.274"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").Value = "'0.9647"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D41").Value = "'0.01693"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.23%  "

$ws.Range("D42").Value = "'6.140"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "'0.4314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.32%  "

$ws.Range("D44").Value = "'104.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.83%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "'7.620"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.36%  "

$ws.Range("D47").Value = "'0.1280"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.94%  "

$ws.Range("D48").Value = "'0.05744"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.18%  "

$ws.Range("D49").Value = "'33.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.73%  "

$ws.Range("D50").Value = "'8.390"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("D51").Value = "'0.3825"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.05%  "
